# T2189_PortfolioValuationCreation.xlsx — "Merged code of Eng"
#
# Content change: AddOpportunity!AA2 was updated from "10.0" to "10000.0".
# Alongside that data edit, the workbook was left with AddOpportunity as the
# active/selected sheet (instead of ValuationPeriod), with a new selected
# cell of AA3 on that sheet.

$wb = $excel.ActiveWorkbook

$wsUsers            = $wb.Worksheets.Item("Users")
$wsValuationPeriod  = $wb.Worksheets.Item("ValuationPeriod")
$wsAddOpportunity   = $wb.Worksheets.Item("AddOpportunity")

# --- the actual data edit -------------------------------------------------
$wsAddOpportunity.Range("AA2").Value = "10000.0"

# --- selections on each sheet (kept/restored as in the saved file) -------
$wsUsers.Range("A2").Select()
$wsValuationPeriod.Range("D2").Select()
$wsAddOpportunity.Range("AA3").Select()

# --- make AddOpportunity the active/selected tab --------------------------
$wsAddOpportunity.Activate()
